$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: fill in Start Time / End Time (Date already present)
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0

# Row 14: new daily record - Date, Start Time, End Time
$ws.Range("A14").Value = 43342
$ws.Range("B14").Value = 0.86111111111111116
$ws.Range("C14").Value = 0.99930555555555556

# Row 15: new daily record - Date, Start Time (End Time not yet recorded)
$ws.Range("A15").Value = 43343
$ws.Range("B15").Value = 0

# Move selection to the next cell to fill in
$ws.Range("C15").Select()
